$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 21407
$ws.Range("I75").Value = 5000
$ws.Range("J75").Value = 26876
$ws.Range("K75").Value = 5000
$ws.Range("L75").Value = 26876
$ws.Range("M75").Value = -4064
$ws.Range("N75").Value = -28748
$ws.Range("H78").Value = 21407
$ws.Range("I78").Value = 5000
$ws.Range("J78").Value = 26876
$ws.Range("K78").Value = 15000
$ws.Range("L78").Value = 80628
$ws.Range("M78").Value = -10320
$ws.Range("N78").Value = -89988
$ws.Range("H116").Value = 4284.8
$ws.Range("I116").Value = 4207.5557
$ws.Range("K116").Value = 4207.5557
$ws.Range("M116").Value = -765.5556999999999
$ws.Range("H137").Value = 6363.303
$ws.Range("I137").Value = 7363.44
$ws.Range("J137").Value = 3237.875
$ws.Range("K137").Value = 22090.32
$ws.Range("L137").Value = 9713.625
$ws.Range("M137").Value = -19540.32
$ws.Range("N137").Value = -14813.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1194108.6
$ws.Range("I32").Value = 1285596.8
$ws.Range("J32").Value = 4762.8
$ws.Range("K32").Value = 1285596.8
$ws.Range("L32").Value = 4762.8
$ws.Range("M32").Value = -1285309.8
$ws.Range("N32").Value = -5336.8
$ws.Range("H74").Value = 1703.22
$ws.Range("I74").Value = 1131.4878
$ws.Range("J74").Value = 4307.778
$ws.Range("K74").Value = 1131.4878
$ws.Range("L74").Value = 4307.778
$ws.Range("M74").Value = -257.4878000000001
$ws.Range("N74").Value = -6055.778
$ws.Range("H77").Value = 1703.22
$ws.Range("I77").Value = 1131.4878
$ws.Range("J77").Value = 4307.778
$ws.Range("K77").Value = 5657.439
$ws.Range("L77").Value = 21538.89
$ws.Range("M77").Value = -1289.439
$ws.Range("N77").Value = -30274.89
$ws.Range("H132").Value = 21417.055
$ws.Range("I132").Value = 31169.03
$ws.Range("J132").Value = 3452.8948
$ws.Range("K132").Value = 93507.09
$ws.Range("L132").Value = 10358.6844
$ws.Range("M132").Value = -90977.09
$ws.Range("N132").Value = -15418.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2516.4412
$ws.Range("I31").Value = 1796.0217
$ws.Range("J31").Value = 4022.7727
$ws.Range("K31").Value = 1796.0217
$ws.Range("L31").Value = 4022.7727
$ws.Range("M31").Value = -1501.0217
$ws.Range("N31").Value = -4612.7727
$ws.Range("H34").Value = 2516.4412
$ws.Range("I34").Value = 1796.0217
$ws.Range("J34").Value = 4022.7727
$ws.Range("K34").Value = 1796.0217
$ws.Range("L34").Value = 4022.7727
$ws.Range("M34").Value = -1594.0217
$ws.Range("N34").Value = -4426.7727
$ws.Range("H58").Value = 1138.2826
$ws.Range("I58").Value = 642.3514
$ws.Range("K58").Value = 642.3514
$ws.Range("M58").Value = -439.3514
$ws.Range("H107").Value = 375.48486
$ws.Range("I107").Value = 325.26086
$ws.Range("J107").Value = 491
$ws.Range("K107").Value = 325.26086
$ws.Range("L107").Value = 491
$ws.Range("M107").Value = 1594.73914
$ws.Range("N107").Value = -4331
$ws.Range("H136").Value = 1138.2826
$ws.Range("I136").Value = 642.3514
$ws.Range("K136").Value = 1927.0542
$ws.Range("M136").Value = 622.9458

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 49
$ws.Range("I12").Value = 17.5
$ws.Range("J12").Value = 70
$ws.Range("K12").Value = 52.5
$ws.Range("L12").Value = 210
$ws.Range("M12").Value = 120.5
$ws.Range("N12").Value = -556
$ws.Range("H64").Value = 1812.25
$ws.Range("I64").Value = 862.3
$ws.Range("J64").Value = 2762.2
$ws.Range("K64").Value = 2586.9
$ws.Range("L64").Value = 8286.599999999999
$ws.Range("M64").Value = -2316.9
$ws.Range("N64").Value = -8826.599999999999
$ws.Range("H67").Value = 1812.25
$ws.Range("I67").Value = 862.3
$ws.Range("J67").Value = 2762.2
$ws.Range("K67").Value = 2586.9
$ws.Range("L67").Value = 8286.599999999999
$ws.Range("M67").Value = -1650.9
$ws.Range("N67").Value = -10158.6
$ws.Range("H87").Value = 5534.8887
$ws.Range("I87").Value = 3478.5
$ws.Range("J87").Value = 7180
$ws.Range("K87").Value = 10435.5
$ws.Range("L87").Value = 21540
$ws.Range("M87").Value = -9187.5
$ws.Range("N87").Value = -24036
$ws.Range("H90").Value = 5534.8887
$ws.Range("I90").Value = 3478.5
$ws.Range("J90").Value = 7180
$ws.Range("K90").Value = 31306.5
$ws.Range("L90").Value = 64620
$ws.Range("M90").Value = -25066.5
$ws.Range("N90").Value = -77100
$ws.Range("H100").Value = 6401.4707
$ws.Range("I100").Value = 4025
$ws.Range("J100").Value = 6550
$ws.Range("K100").Value = 12075
$ws.Range("L100").Value = 19650
$ws.Range("M100").Value = -11264
$ws.Range("N100").Value = -21272
$ws.Range("H102").Value = 7375
$ws.Range("J102").Value = 7428.5713
$ws.Range("L102").Value = 22285.7139
$ws.Range("N102").Value = -27153.7139
$ws.Range("H114").Value = 2425.8125
$ws.Range("I114").Value = 4320.2856
$ws.Range("J114").Value = 952.3333
$ws.Range("K114").Value = 12960.8568
$ws.Range("L114").Value = 2856.9999
$ws.Range("M114").Value = -9706.856800000001
$ws.Range("N114").Value = -9364.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3239.4146
$ws.Range("I126").Value = 3128.087
$ws.Range("J126").Value = 3381.6667
$ws.Range("K126").Value = 9384.261
$ws.Range("L126").Value = 10145.0001
$ws.Range("M126").Value = -6914.261
$ws.Range("N126").Value = -15085.0001
$ws.Range("H132").Value = 3748.9302
$ws.Range("I132").Value = 3617.8235
$ws.Range("J132").Value = 4244.222
$ws.Range("K132").Value = 10853.4705
$ws.Range("L132").Value = 12732.666
$ws.Range("M132").Value = -8323.470499999999
$ws.Range("N132").Value = -17792.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1081.9445
$ws.Range("I16").Value = 950.24243
$ws.Range("J16").Value = 2530.6667
$ws.Range("K16").Value = 950.24243
$ws.Range("L16").Value = 2530.6667
$ws.Range("M16").Value = -780.24243
$ws.Range("N16").Value = -2870.6667
$ws.Range("H61").Value = 1461.375
$ws.Range("I61").Value = 892.7
$ws.Range("K61").Value = 892.7
$ws.Range("M61").Value = -690.7
$ws.Range("H113").Value = 1461.375
$ws.Range("I113").Value = 892.7
$ws.Range("K113").Value = 892.7
$ws.Range("M113").Value = 1277.3
$ws.Range("H122").Value = 2356.2
$ws.Range("I122").Value = 2440.7144
$ws.Range("J122").Value = 2248.6365
$ws.Range("K122").Value = 7322.1432
$ws.Range("L122").Value = 6745.9095
$ws.Range("M122").Value = -4872.1432
$ws.Range("N122").Value = -11645.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1988.8
$ws.Range("I96").Value = 2033.3334
$ws.Range("J96").Value = 1969.7142
$ws.Range("K96").Value = 2033.3334
$ws.Range("L96").Value = 1969.7142
$ws.Range("M96").Value = -660.3334
$ws.Range("N96").Value = -4715.7142
$ws.Range("H122").Value = 1568.5333
$ws.Range("I122").Value = 1574.5294
$ws.Range("K122").Value = 4723.5882
$ws.Range("M122").Value = -2273.5882
$ws.Range("H126").Value = 1276.3914
$ws.Range("I126").Value = 657.25
$ws.Range("K126").Value = 1971.75
$ws.Range("M126").Value = 498.25
